{"js": "// Update each \"AxB=\" multiplication expression in the document's table\n// to its new value. Pairs are listed in document order so that the\n// one value that is both an old and a new string (\"91\u00d742=\") is handled\n// unambiguously: the occurrence that must change away from \"91\u00d742=\"\n// is processed (and therefore no longer matches) before the later\n// cell is changed to become \"91\u00d742=\".\nconst replacements = [\n  [\"75\u00d778=\", \"47\u00d790=\"],\n  [\"41\u00d711=\", \"57\u00d761=\"],\n  [\"95\u00d716=\", \"59\u00d767=\"],\n  [\"54\u00d757=\", \"78\u00d783=\"],\n  [\"91\u00d742=\", \"40\u00d773=\"],\n  [\"61\u00d728=\", \"53\u00d776=\"],\n  [\"46\u00d743=\", \"31\u00d785=\"],\n  [\"67\u00d757=\", \"82\u00d797=\"],\n  [\"15\u00d718=\", \"40\u00d713=\"],\n  [\"32\u00d735=\", \"54\u00d777=\"],\n  [\"66\u00d724=\", \"91\u00d742=\"],\n  [\"29\u00d739=\", \"39\u00d752=\"],\n  [\"62\u00d763=\", \"50\u00d715=\"],\n  [\"56\u00d789=\", \"49\u00d799=\"],\n  [\"37\u00d715=\", \"36\u00d724=\"],\n  [\"65\u00d743=\", \"46\u00d799=\"],\n  [\"30\u00d770=\", \"50\u00d778=\"],\n  [\"94\u00d773=\", \"98\u00d747=\"],\n  [\"69\u00d714=\", \"67\u00d744=\"],\n  [\"20\u00d745=\", \"81\u00d799=\"],\n  [\"77\u00d756=\", \"41\u00d717=\"],\n  [\"18\u00d785=\", \"29\u00d748=\"],\n  [\"52\u00d744=\", \"12\u00d713=\"],\n  [\"50\u00d767=\", \"85\u00d714=\"],\n  [\"90\u00d794=\", \"81\u00d789=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  // Replace only the first match; each source string is unique at the\n  // moment we search for it given the document-order processing above.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Update each \"AxB=\" multiplication expression in the document's table\n# to its new value. Pairs are listed in document order so that the one\n# value that is both an old and a new string (\"91\u00d742=\") is handled\n# unambiguously: the occurrence that must change away from \"91\u00d742=\"\n# is processed (and therefore no longer matches) before the later cell\n# is changed to become \"91\u00d742=\".\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"75\u00d778=\", \"47\u00d790=\"),\n    @(\"41\u00d711=\", \"57\u00d761=\"),\n    @(\"95\u00d716=\", \"59\u00d767=\"),\n    @(\"54\u00d757=\", \"78\u00d783=\"),\n    @(\"91\u00d742=\", \"40\u00d773=\"),\n    @(\"61\u00d728=\", \"53\u00d776=\"),\n    @(\"46\u00d743=\", \"31\u00d785=\"),\n    @(\"67\u00d757=\", \"82\u00d797=\"),\n    @(\"15\u00d718=\", \"40\u00d713=\"),\n    @(\"32\u00d735=\", \"54\u00d777=\"),\n    @(\"66\u00d724=\", \"91\u00d742=\"),\n    @(\"29\u00d739=\", \"39\u00d752=\"),\n    @(\"62\u00d763=\", \"50\u00d715=\"),\n    @(\"56\u00d789=\", \"49\u00d799=\"),\n    @(\"37\u00d715=\", \"36\u00d724=\"),\n    @(\"65\u00d743=\", \"46\u00d799=\"),\n    @(\"30\u00d770=\", \"50\u00d778=\"),\n    @(\"94\u00d773=\", \"98\u00d747=\"),\n    @(\"69\u00d714=\", \"67\u00d744=\"),\n    @(\"20\u00d745=\", \"81\u00d799=\"),\n    @(\"77\u00d756=\", \"41\u00d717=\"),\n    @(\"18\u00d785=\", \"29\u00d748=\"),\n    @(\"52\u00d744=\", \"12\u00d713=\"),\n    @(\"50\u00d767=\", \"85\u00d714=\"),\n    @(\"90\u00d794=\", \"81\u00d789=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    # wrap = wdFindContinue (1); replace = wdReplaceOne (1) so only the\n    # first (next) matching occurrence is changed.\n    $found = $d.Content.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 1)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
